$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.01925380265402739
$ws.Cells.Item(2, 3).Value = 2888070398.104109
$ws.Cells.Item(2, 4).Value = 12342242302.91712
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "2025-03-09T04:26:47.320735"

$ws.Cells.Item(3, 2).Value = 0.0000004238788545257737
$ws.Cells.Item(3, 3).Value = 17802911.89008249
$ws.Cells.Item(3, 4).Value = 46599258.27682871
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "2025-03-09T04:26:47.320735"

$ws.Cells.Item(4, 2).Value = 0.004240403655724095
$ws.Cells.Item(4, 3).Value = 2497597753221.492
$ws.Cells.Item(4, 4).Value = 4227602636016.839
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "2025-03-09T04:26:47.320735"

$ws.Cells.Item(5, 2).Value = 1.112345263355795
$ws.Cells.Item(5, 3).Value = 111234526335.5795
$ws.Cells.Item(5, 4).Value = 168272317101.9254
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "2025-03-09T04:26:47.320735"

$ws.Cells.Item(6, 2).Value = 228.1441865290555
$ws.Cells.Item(6, 3).Value = 3422162797935.833
$ws.Cells.Item(6, 4).Value = 3544438988512.615
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "2025-03-09T04:26:47.320735"

$ws.Cells.Item(7, 2).Value = 378.7045994700504
$ws.Cells.Item(7, 3).Value = 2809988128067.774
$ws.Cells.Item(7, 4).Value = 3565010145491.38
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "2025-03-09T04:26:47.320735"
